# APSS BOM Template update
#
# Commit: "Updating output job file and BOM template" — the BOM template
# portion of that commit retargets the merge/report field placeholders in
# column E (rows 4-6) so they refer to the correct parameter names exported
# by the PCB tool (PCBName -> DesignName, PCBRevision -> BoardRevision,
# PCBDesigner -> Designer), and leaves the cursor/selection resting on F5
# (one row up from where it previously sat on F6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "PCB Name:"      | Field=PCBName      -> Field=DesignName
# Row 5: "PCB Revision:"  | Field=PCBRevision  -> Field=BoardRevision
# Row 6: "PCB Designer:"  | Field=PCBDesigner  -> Field=Designer
$ws.Range("E4").Value = "Field=DesignName"
$ws.Range("E5").Value = "Field=BoardRevision"
$ws.Range("E6").Value = "Field=Designer"

# Leave the active selection on F5 (matches the saved view state in the diff).
$ws.Range("F5").Select()
